# Updated cryptos list on Sun Jun 30 09:18:58 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "61.549.79";  E = "  +1.07%  " },
    @{ Row = 3;  D = "3.393.40";   E = "  +0.09%  " },
    @{ Row = 4;  D = $null;        E = "  +0.00%  " },
    @{ Row = 5;  D = "576.79";     E = "  +0.99%  " },
    @{ Row = 6;  D = "141.30";     E = "  -0.57%  " },
    @{ Row = 7;  D = $null;        E = "  +0.08%  " },
    @{ Row = 8;  D = $null;        E = "  -0.20%  " },
    @{ Row = 9;  D = $null;        E = "  +2.40%  " },
    @{ Row = 10; D = $null;        E = "  -0.58%  " },
    @{ Row = 11; D = $null;        E = "  -1.76%  " },
    @{ Row = 12; D = "3.974.51";   E = "  +0.15%  " },
    @{ Row = 13; D = $null;        E = "  +0.33%  " },
    @{ Row = 14; D = $null;        E = "  +0.34%  " },
    @{ Row = 15; D = "3.391.57";   E = "  -0.05%  " },
    @{ Row = 16; D = $null;        E = "  +0.05%  " },
    @{ Row = 17; D = "61.563.07";  E = "  +0.95%  " },
    @{ Row = 18; D = "6.14";       E = "  -0.50%  " },
    @{ Row = 19; D = "13.70";      E = "  -1.18%  " },
    @{ Row = 20; D = "8.97";       E = "  -0.04%  " },
    @{ Row = 21; D = "391.11";     E = "  +1.70%  " },
    @{ Row = 22; D = "75.52";      E = "  +1.49%  " },
    @{ Row = 23; D = "0.554";      E = "  -0.75%  " },
    @{ Row = 24; D = $null;        E = "  +0.01%  " },
    @{ Row = 25; D = "0.0000113";  E = "  -3.92%  " },
    @{ Row = 26; D = $null;        E = "  +8.74%  " },
    @{ Row = 27; D = $null;        E = "  -0.03%  " },
    @{ Row = 28; D = "7.28";       E = "  -1.45%  " },
    @{ Row = 29; D = "8.04";       E = "  +0.74%  " },
    @{ Row = 30; D = $null;        E = "  +0.40%  " },
    @{ Row = 31; D = $null;        E = "  -0.07%  " },
    @{ Row = 32; D = $null;        E = "  -4.37%  " },
    @{ Row = 33; D = "23.40";      E = "  -0.44%  " },
    @{ Row = 34; D = $null;        E = "  -0.54%  " },
    @{ Row = 35; D = "167.39";     E = "  +0.21%  " },
    @{ Row = 36; D = "5.04";       E = "  +1.43%  " },
    @{ Row = 37; D = "3.428.24";   E = "  +0.22%  " },
    @{ Row = 38; D = $null;        E = "  -0.66%  " },
    @{ Row = 39; D = "0.0773";     E = "  -0.20%  " },
    @{ Row = 40; D = "26.10";      E = "  -4.13%  " },
    @{ Row = 41; D = $null;        E = "  +0.00%  " },
    @{ Row = 42; D = $null;        E = "  +0.14%  " },
    @{ Row = 43; D = $null;        E = "  -0.58%  " },
    @{ Row = 44; D = $null;        E = "  +1.03%  " },
    @{ Row = 45; D = "2.468.55";   E = "  -0.40%  " },
    @{ Row = 46; D = "23.09";      E = "  +0.18%  " },
    @{ Row = 47; D = "6.69";       E = "  -1.96%  " },
    @{ Row = 48; D = "1.00";       E = "  +0.03%  " },
    @{ Row = 49; D = $null;        E = "  -1.92%  " },
    @{ Row = 50; D = "2.07";       E = "  -0.93%  " },
    @{ Row = 51; D = $null;        E = "  -1.73%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        # Column D holds price strings like "61.549.79" that look numeric to
        # Excel's auto-detection. Force text storage (matching the original
        # t="inlineStr"/text cells) by temporarily applying a text number
        # format, then restore the default "Normal" style so no stray
        # formatting is left behind on the cell.
        $cellD = $ws.Range("D$r")
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }

    if ($null -ne $u.E) {
        # Column E values (e.g. "  +1.07%  ") are never numeric-looking
        # because of the surrounding whitespace, so a plain value
        # assignment keeps them as text.
        $ws.Range("E$r").Value = $u.E
    }
}
